$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that sits in the title paragraph
#    ("One Week at a Time Data Model").
try {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldBm.Delete()
} catch {
    # bookmark may not exist; ignore
}

# 2. Find the "Name of task category" paragraph (inside the table) and the
#    "Unique" paragraph that immediately follows it, then remove the
#    "Unique" paragraph entirely.
$targetPara = $null
$uniquePara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "Name of task category") {
        $targetPara = $p
    } elseif (($targetPara -ne $null) -and ($uniquePara -eq $null) -and ($t -eq "Unique")) {
        $uniquePara = $p
    }
}

if ($uniquePara -ne $null) {
    $uniquePara.Range.Delete()
}

# 3. Re-create the "_GoBack" bookmark, now positioned at the end of the
#    "Name of task category" paragraph content (immediately before the
#    paragraph mark), matching the target layout.
#    Collapsed (zero-length) ranges aren't placed correctly by
#    Bookmarks.Add in this environment, so build the insertion point by
#    temporarily inserting a placeholder character, wrapping the bookmark
#    around it, and then deleting the placeholder -- leaving a
#    zero-length bookmark exactly where the placeholder was.
$r = $targetPara.Range
$null = $r.MoveEnd(1, -1) # wdCharacter: exclude the trailing paragraph mark
$null = $r.Collapse(0)    # wdCollapseEnd: collapse to just after "category"
$insertStart = $r.Start
$null = $r.InsertAfter([char]0x0001)
$null = $d.Bookmarks.Add("_GoBack", $r)
$placeholder = $d.Range($insertStart, $insertStart + 1)
$null = $placeholder.Delete()
